$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.572.36"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.603.56"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.489"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "4.215.96"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "3.609.93"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "66.654.00"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "428.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "3.753.35"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.09%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").Value = "3.601.19"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("E33").Value = "  +4.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "178.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.957"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "2.433.00"
$ws.Range("E51").Value = "  +5.74%  "
